# "Error Calculations and Plots"
# Two source rows (ID "RM 232" and ID "SC 92") were dropped from the
# missing-data table, shifting every row below them up, and the
# imputed/missing-value pattern in columns C, E and F was re-rolled for
# several of the remaining rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the two rows -----------------------------------------------
# Row 26 held ID "RM 232"; after removing it, the row that used to be 28
# ("SC 92") is now row 27. Deleting shifts everything below up by one
# each time.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# --- Re-roll the missing-value pattern on the remaining rows ---------
$ws.Range("E3").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("F8").Value = 17.05
$ws.Range("F10").Value = 16.43
$ws.Range("F12").Value = ""
$ws.Range("F15").Value = 16.2
$ws.Range("F18").Value = ""
$ws.Range("F19").Value = ""
$ws.Range("F25").Value = 16.6
$ws.Range("C26").Value = 10.8
$ws.Range("C27").Value = ""
$ws.Range("F29").Value = ""
$ws.Range("C33").Value = 10.4
$ws.Range("E33").Value = -10.7
